# Refresh the cryptocurrency price/volume table (cryptos.xlsx)
# as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while preserving it as literal text,
# even when the text happens to look like a number (e.g. '69.47').
# Excel auto-converts plain numeric-looking strings assigned via .Value
# into real floating point numbers, which would corrupt the data (and
# introduce binary floating point noise, e.g. 69.47 -> 69.4699999...).
# Temporarily forcing a Text number format avoids that, and resetting
# the style back to Normal afterwards leaves the cell without any
# custom formatting, just like the source file.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '69.555.75'
$ws.Range("E2").Value = '  +0.44%  '
# Row 3
$ws.Range("D3").Value = '2.494.10'
$ws.Range("E3").Value = '  -0.99%  '
# Row 5
Set-TextValue $ws.Range("D5") '570.38'
$ws.Range("E5").Value = '  -0.35%  '
# Row 6
Set-TextValue $ws.Range("D6") '166.68'
$ws.Range("E6").Value = '  +0.45%  '
# Row 7
$ws.Range("E7").Value = '  -0.05%  '
# Row 8
$ws.Range("E8").Value = '  -1.20%  '
# Row 9
$ws.Range("E9").Value = '  -0.25%  '
# Row 10
$ws.Range("E10").Value = '  -0.61%  '
# Row 11
$ws.Range("E11").Value = '  -1.52%  '
# Row 12
Set-TextValue $ws.Range("D12") '4.87'
$ws.Range("E12").Value = '  -0.64%  '
# Row 13
$ws.Range("D13").Value = '2.951.35'
$ws.Range("E13").Value = '  -1.04%  '
# Row 14
$ws.Range("D14").Value = '69.459.47'
$ws.Range("E14").Value = '  +0.53%  '
# Row 15
$ws.Range("E15").Value = '  -0.38%  '
# Row 16
Set-TextValue $ws.Range("D16") '24.25'
$ws.Range("E16").Value = '  -2.30%  '
# Row 17
$ws.Range("D17").Value = '2.499.15'
$ws.Range("E17").Value = '  -0.84%  '
# Row 18
Set-TextValue $ws.Range("D18") '11.26'
$ws.Range("E18").Value = '  -0.70%  '
# Row 19
Set-TextValue $ws.Range("D19") '354.72'
$ws.Range("E19").Value = '  +1.89%  '
# Row 20
Set-TextValue $ws.Range("D20") '7.39'
$ws.Range("E20").Value = '  -3.23%  '
# Row 21
$ws.Range("E21").Value = '  -0.49%  '
# Row 22
$ws.Range("E22").Value = '  -4.78%  '
# Row 23
$ws.Range("E23").Value = '  -0.13%  '
# Row 24
Set-TextValue $ws.Range("D24") '69.47'
$ws.Range("E24").Value = '  -1.35%  '
# Row 25
Set-TextValue $ws.Range("D25") '3.82'
$ws.Range("E25").Value = '  -3.62%  '
# Row 26
$ws.Range("D26").Value = '2.623.76'
$ws.Range("E26").Value = '  -0.79%  '
# Row 27
Set-TextValue $ws.Range("D27") '8.62'
$ws.Range("E27").Value = '  -3.00%  '
# Row 28
Set-TextValue $ws.Range("D28") '0.998'
$ws.Range("E28").Value = '  -0.06%  '
# Row 29
$ws.Range("E29").Value = '  -2.13%  '
# Row 30
$ws.Range("E30").Value = '  -2.44%  '
# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D31") '440.03'
$ws.Range("E31").Value = '  -4.86%  '
# Row 32
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D32") '1.19'
$ws.Range("E32").Value = '  -3.39%  '
# Row 33
Set-TextValue $ws.Range("D33") '1.00'
$ws.Range("E33").Value = '  +0.00%  '
# Row 34
$ws.Range("E34").Value = '  -1.06%  '
# Row 35
$ws.Range("B35").Value = 'POPCAT'
$ws.Range("C35").Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
Set-TextValue $ws.Range("D35") '2.90'
$ws.Range("E35").Value = '  +84.43%  '
# Row 36
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D36") '154.02'
$ws.Range("E36").Value = '  -2.20%  '
# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D37") '0.113'
$ws.Range("E37").Value = '  -2.96%  '
# Row 38
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D38") '19.07'
$ws.Range("E38").Value = '  +0.33%  '
# Row 39
$ws.Range("E39").Value = '  -1.92%  '
# Row 40
$ws.Range("E40").Value = '  +0.03%  '
# Row 41
$ws.Range("E41").Value = '  -1.51%  '
# Row 42
$ws.Range("E42").Value = '  -2.29%  '
# Row 43
$ws.Range("E43").Value = '  -1.77%  '
# Row 44
Set-TextValue $ws.Range("D44") '2.18'
$ws.Range("E44").Value = '  -3.10%  '
# Row 45
$ws.Range("E45").Value = '  -5.22%  '
# Row 46
Set-TextValue $ws.Range("D46") '138.68'
$ws.Range("E46").Value = '  -2.19%  '
# Row 47
$ws.Range("E47").Value = '  -1.15%  '
# Row 48
Set-TextValue $ws.Range("D48") '0.506'
$ws.Range("E48").Value = '  -3.05%  '
# Row 49
Set-TextValue $ws.Range("D49") '0.0723'
$ws.Range("E49").Value = '  -1.00%  '
# Row 50
$ws.Range("E50").Value = '  -1.49%  '
# Row 51
Set-TextValue $ws.Range("D51") '0.0924'
$ws.Range("E51").Value = '  -0.87%  '
